# Automatische test-sync: 2025-06-30 19:44:50
# Adds a new "Testmail #4" row (row 5) to the "Logs" sheet, extends the
# conditional-formatting ranges that covered rows 2-4 so they also cover
# row 5, and re-syncs the "Dashboard" category/count table (which is kept
# sorted by descending count) to reflect the new totals.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- 1. Append the new log entry in row 5 of the "Logs" sheet -------------
$logs.Range("A5").Value = "Ik stuur het pakket terug."
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("C5").Value = "Testmail #4: Ik stuur het pakket terug."
$logs.Range("D5").Value = "Retour / Terugbetaling"
$logs.Range("E5").Value = "Beste klant,`nBedankt voor uw bericht. Om uw retourzending zo soepel mogelijk te laten verlopen, vragen wij u vriendelijk om het volgende te doen:`n- Vul het retourformulier in dat bij uw bestelling zat en voeg dit toe aan het pakket.`n- Stuur het pakket terug naar het volgende adres: [adres retourzending].`n- Zodra wij uw retourzending hebben ontvangen, zullen wij het verder afhandelen en u op de hoogte houden van de status van uw retour.`nMocht u nog verdere vragen of opmerkingen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Range("F5").Value = "2025-06-30 19:44:16"
$logs.Range("G5").Value = "Ja"
$logs.Range("H5").Value = "Nee"
$logs.Range("I5").Value = "Ja"
$logs.Range("J5").Value = "Nee"

# --- 2. Extend the conditional-formatting blocks that applied to rows 2-4
#        (D/G/H/I/J) so that they also cover the newly added row 5. Each
#        block shares a single sqref, so modifying any one rule in a block
#        moves the whole block. -------------------------------------------
$logs.Range("D2:D4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D5"))
$logs.Range("G2:G4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G5"))
$logs.Range("H2:H4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H5"))
$logs.Range("I2:I4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I5"))
$logs.Range("J2:J4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J5"))

# --- 3. Re-sync the "Dashboard" category summary table ---------------------
# It stays sorted by descending count; "Retour / Terugbetaling" now has 2
# occurrences (rows 4 and 5 in "Logs") and moves to the top, pushing the
# other two categories down one row.
$dash.Range("A2").Value = "Retour / Terugbetaling"
$dash.Range("B2").Value = 2
$dash.Range("A3").Value = "Openingstijden / Locatie"
$dash.Range("B3").Value = 1
$dash.Range("A4").Value = "Bestelling / Levering"
$dash.Range("B4").Value = 1
